# Weekly price update: insert a new daily record right before the existing
# row for this series (old row 219), pushing all subsequent rows down by one.
# The new row carries the latest reading (Fecha 45001) and the remaining
# rows keep their prior values, just shifted down one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 219; everything from old row 219 onward shifts
# down to row 220 onward (row 219 inherits formatting from the row above,
# same as Excel's native "Insert Row" behavior).
$ws.Rows(219).Insert()

$ws.Range("A219").Value = 4
$ws.Range("B219").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C219").Value = "Los Lagos"
$ws.Range("D219").Value = 45001
$ws.Range("E219").Value = 10
$ws.Range("F219").Value = 100112039
$ws.Range("G219").Value = "Ciboulette"
$ws.Range("H219").Value = "Sin especificar"
$ws.Range("I219").Value = "Primera"
$ws.Range("J219").Value = 120
$ws.Range("K219").Value = 3500
$ws.Range("L219").Value = 3500
$ws.Range("M219").Value = 3500
$ws.Range("N219").Value = "`$/docena de atados"
$ws.Range("O219").Value = "Región Metropolitana"
$ws.Range("P219").Value = 1167
$ws.Range("Q219").Value = 3
$ws.Range("R219").Value = "Hortaliza"
